# Apply crypto price/volume updates (GitHub Actions daily refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "46.054.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.614.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.596"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.83%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.580"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.83"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0841"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.02"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.10"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.003.95"
$ws.Range("D14").Style = "Normal"
$ws.Range("E15").Value = "  +0.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.608.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.913"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.82"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "46.056.87"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.72%  "
$ws.Range("E20").Value = "  -1.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.73"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.43%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "74.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.49%  "
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "285.55"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +11.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "29.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.04%  "
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("E29").Value = "  +0.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.57"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "38.78"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.21%  "
$ws.Range("E32").Value = "  -2.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.26"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.59%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.29"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.65%  "
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "157.13"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.98%  "
$ws.Range("E37").Value = "  -2.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0833"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.122"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.65%  "
$ws.Range("E40").Value = "  +0.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "15.97"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.50%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0328"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.05"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.58"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.55"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.109.57"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.999"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "94.20"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.21"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.72%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "109.21"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.49%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.200"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.50%  "
